$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-11 Monday" "2024-11-12 Tuesday"

Replace-Text "302×4=" "815×9="
Replace-Text "908×8=" "431×3="
Replace-Text "214×3=" "463×4="
Replace-Text "141×8=" "946×4="
Replace-Text "359×8=" "743×5="
Replace-Text "658×7=" "340×3="
Replace-Text "337×9=" "439×5="
Replace-Text "394×9=" "220×6="
Replace-Text "434×2=" "705×3="
Replace-Text "401×8=" "288×3="
Replace-Text "106×6=" "149×6="
Replace-Text "264×9=" "300×8="
Replace-Text "198×6=" "630×8="
Replace-Text "336×5=" "670×7="
Replace-Text "832×2=" "165×9="
Replace-Text "758×5=" "928×5="
Replace-Text "514×6=" "406×5="
Replace-Text "551×5=" "912×5="
Replace-Text "960×5=" "510×2="
Replace-Text "937×4=" "685×8="
Replace-Text "142×6=" "165×7="
Replace-Text "236×7=" "262×9="
Replace-Text "767×5=" "865×8="
Replace-Text "774×6=" "842×3="
Replace-Text "443×9=" "152×3="
